$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2661
$ws.Range("J17").Value = 2773.652
$ws.Range("L17").Value = 8320.956
$ws.Range("N17").Value = -8656.956
$ws.Range("H62").Value = 3432.182
$ws.Range("I62").Value = 3432.182
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3432.182
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2808.182
$ws.Range("H65").Value = 3432.182
$ws.Range("I65").Value = 3432.182
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17160.91
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -14040.91
$ws.Range("H132").Value = 7658.5835
$ws.Range("I132").Value = 5308.0938
$ws.Range("J132").Value = 26462.5
$ws.Range("K132").Value = 15924.2814
$ws.Range("L132").Value = 79387.5
$ws.Range("M132").Value = -13394.2814
$ws.Range("N132").Value = -84447.5
$ws.Range("H137").Value = 2611.28
$ws.Range("I137").Value = 3628.8
$ws.Range("J137").Value = 1932.9333
$ws.Range("K137").Value = 10886.4
$ws.Range("L137").Value = 5798.7999
$ws.Range("M137").Value = -8336.400000000001
$ws.Range("N137").Value = -10898.7999
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1885.16
$ws.Range("I2").Value = 1135.2354
$ws.Range("J2").Value = 3478.75
$ws.Range("K2").Value = 1135.2354
$ws.Range("L2").Value = 3478.75
$ws.Range("M2").Value = -1022.2354
$ws.Range("N2").Value = -3704.75
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("H116").Value = 1885.16
$ws.Range("I116").Value = 1135.2354
$ws.Range("J116").Value = 3478.75
$ws.Range("K116").Value = 1135.2354
$ws.Range("L116").Value = 3478.75
$ws.Range("M116").Value = 1158.7646
$ws.Range("N116").Value = -8066.75
$ws.Range("H122").Value = 2309.3
$ws.Range("I122").Value = 1664.5769
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 4993.7307
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -2543.7307
$ws.Range("N122").Value = -24400
$ws.Range("H132").Value = 4292.1177
$ws.Range("I132").Value = 2348.476
$ws.Range("J132").Value = 5652.6665
$ws.Range("K132").Value = 7045.428
$ws.Range("L132").Value = 16957.9995
$ws.Range("M132").Value = -4515.428
$ws.Range("N132").Value = -22017.9995
$ws.Range("N104").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1885.16
$ws.Range("I3").Value = 1135.2354
$ws.Range("J3").Value = 3478.75
$ws.Range("K3").Value = 1135.2354
$ws.Range("L3").Value = 3478.75
$ws.Range("M3").Value = -1021.2354
$ws.Range("N3").Value = -3706.75
$ws.Range("H20").Value = 3445
$ws.Range("I20").Value = 2694.6667
$ws.Range("J20").Value = 3945.2222
$ws.Range("K20").Value = 2694.6667
$ws.Range("L20").Value = 3945.2222
$ws.Range("M20").Value = -2447.6667
$ws.Range("N20").Value = -4439.2222
$ws.Range("H94").Value = 401.35
$ws.Range("I94").Value = 353.94116
$ws.Range("K94").Value = 353.94116
$ws.Range("M94").Value = 97.05883999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2393.6785
$ws.Range("I16").Value = 1900.9584
$ws.Range("J16").Value = 5350
$ws.Range("K16").Value = 1900.9584
$ws.Range("L16").Value = 5350
$ws.Range("M16").Value = -1613.9584
$ws.Range("N16").Value = -5924
$ws.Range("H22").Value = 873.6667
$ws.Range("I22").Value = 745.6
$ws.Range("J22").Value = 1033.75
$ws.Range("K22").Value = 745.6
$ws.Range("L22").Value = 1033.75
$ws.Range("M22").Value = -395.6
$ws.Range("N22").Value = -1733.75
$ws.Range("H105").Value = 2108.3333
$ws.Range("I105").Value = 1716.6666
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1716.6666
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 30.33339999999998
$ws.Range("N105").Value = -5994
$ws.Range("H113").Value = 2393.6785
$ws.Range("I113").Value = 1900.9584
$ws.Range("J113").Value = 5350
$ws.Range("K113").Value = 1900.9584
$ws.Range("L113").Value = 5350
$ws.Range("M113").Value = 269.0416
$ws.Range("N113").Value = -9690
$ws.Range("H122").Value = 3003
$ws.Range("I122").Value = 1188.8889
$ws.Range("J122").Value = 4487.273
$ws.Range("K122").Value = 3566.6667
$ws.Range("L122").Value = 13461.819
$ws.Range("M122").Value = -1116.6667
$ws.Range("N122").Value = -18361.819
$ws.Range("H141").Value = 44945
$ws.Range("J141").Value = 44945
$ws.Range("L141").Value = 44945
$ws.Range("N141").Value = -55305

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1373.48
$ws.Range("I113").Value = 1329.3158
$ws.Range("J113").Value = 1513.3334
$ws.Range("K113").Value = 1329.3158
$ws.Range("L113").Value = 1513.3334
$ws.Range("M113").Value = 840.6841999999999
$ws.Range("N113").Value = -5853.3334
$ws.Range("H132").Value = 2621.4
$ws.Range("J132").Value = 2655.3333
$ws.Range("L132").Value = 7965.999899999999
$ws.Range("N132").Value = -13025.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 522.6896400000001
$ws.Range("I22").Value = 286.07693
$ws.Range("J22").Value = 714.9375
$ws.Range("K22").Value = 286.07693
$ws.Range("L22").Value = 714.9375
$ws.Range("M22").Value = 8.923069999999996
$ws.Range("N22").Value = -1304.9375
$ws.Range("H27").Value = 522.6896400000001
$ws.Range("I27").Value = 286.07693
$ws.Range("J27").Value = 714.9375
$ws.Range("K27").Value = 286.07693
$ws.Range("L27").Value = 714.9375
$ws.Range("M27").Value = -179.07693
$ws.Range("N27").Value = -928.9375
$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -14992
$ws.Range("H122").Value = 6144.8667
$ws.Range("I122").Value = 7253.3184
$ws.Range("J122").Value = 3096.625
$ws.Range("K122").Value = 21759.9552
$ws.Range("L122").Value = 9289.875
$ws.Range("M122").Value = -19309.9552
$ws.Range("N122").Value = -14189.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7007.231
$ws.Range("I96").Value = 3042.3333
$ws.Range("J96").Value = 10405.714
$ws.Range("K96").Value = 3042.3333
$ws.Range("L96").Value = 10405.714
$ws.Range("M96").Value = -1669.3333
$ws.Range("N96").Value = -13151.714
$ws.Range("H100").Value = 613.4286
$ws.Range("I100").Value = 545.2308
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1090.4616
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -549.4616000000001
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 5922.769
$ws.Range("I122").Value = 7221.8096
$ws.Range("J122").Value = 4407.222
$ws.Range("K122").Value = 21665.4288
$ws.Range("L122").Value = 13221.666
$ws.Range("M122").Value = -19215.4288
$ws.Range("N122").Value = -18121.666
$ws.Range("H126").Value = 5862.2856
$ws.Range("I126").Value = 6788.6523
$ws.Range("J126").Value = 1601
$ws.Range("K126").Value = 20365.9569
$ws.Range("L126").Value = 4803
$ws.Range("M126").Value = -17895.9569
$ws.Range("N126").Value = -9743
